# Adapt column header formatting to respective input file names (FV2410 / FV2504)
# and turn the sheet's data range into a real Excel Table, with the header
# row frozen in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Rename the header cells: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

# Columns A..J hold the "_old" -> "_FV2410" headers, column K holds "diff"
# (unchanged), columns L..U hold the "_new" -> "_FV2504" headers.
for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $fv2410Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $fv2504Headers[$i]
}

# --- 2. Turn A1:U89 into a real Excel Table ---
$rng = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
